$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117; this pushes the existing rows 117-175
# down to 118-176 (matching the diff, where every row's data shifts down by
# one and a brand-new row of "Arándano (blue)" data is inserted at the top
# of this block).
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new data record.
$ws.Range("A117").Value = 9
$ws.Range("B117").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C117").Value = "Metropolitana"
$ws.Range("D117").Value = 44609
$ws.Range("E117").Value = 13
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100101
$ws.Range("H117").Value = "Berries"
$ws.Range("I117").Value = 100101001
$ws.Range("J117").Value = "Arándano (blue)"
$ws.Range("K117").Value = "Sin especificar"
$ws.Range("L117").Value = "Primera"
$ws.Range("M117").Value = 440
$ws.Range("N117").Value = 3600
$ws.Range("O117").Value = 3600
$ws.Range("P117").Value = 3600
$ws.Range("Q117").Value = "$/bandeja 2 kilos"
$ws.Range("R117").Value = "Provincia de Curicó"
$ws.Range("S117").Value = 1800
$ws.Range("T117").Value = 2
